# aggiornamento a l 23 agosto 2021
# Append daily rows for 2021-08-10 .. 2021-08-23 (Excel serials 44418..44431)
# to the bottom of the single data table on Sheet1 (columns A:D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44418, 1, 5, 109.051254089422),
    @(44419, 0, 4, 87.24100327153762),
    @(44420, 1, 5, 109.051254089422),
    @(44421, 1, 4, 87.24100327153762),
    @(44422, 1, 5, 109.051254089422),
    @(44423, 0, 5, 109.051254089422),
    @(44424, 2, 6, 130.8615049073064),
    @(44425, 0, 5, 109.051254089422),
    @(44426, 0, 5, 109.051254089422),
    @(44427, 0, 4, 87.24100327153762),
    @(44428, 0, 3, 65.43075245365321),
    @(44429, 0, 2, 43.62050163576881),
    @(44430, 0, 2, 43.62050163576881),
    @(44431, 0, 0, 0)
)

# Last existing populated row is 343 (date serial 44417); new rows continue below it.
$lastRow = 343
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Copy column A's cell format (date style with border/bold/alignment) down
    # from the last existing row, then overwrite with the new date value so
    # the style index is reused rather than duplicated.
    $ws.Cells.Item($lastRow, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
